$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.441.74'

$ws.Range("D3").Value = '1.827.44'
$ws.Range("E3").Value = '  +9.08%  '

$ws.Range("D5").Value = '230.44'
$ws.Range("E5").Value = '  +4.91%  '

$ws.Range("D6").Value = '0.575'
$ws.Range("E6").Value = '  +8.85%  '

$ws.Range("E7").Value = '  -0.12%  '

$ws.Range("D8").Value = '31.56'
$ws.Range("E8").Value = '  +7.77%  '

$ws.Range("D9").Value = '46.84'
$ws.Range("E9").Value = '  +5.76%  '

$ws.Range("D10").Value = '0.289'
$ws.Range("E10").Value = '  +9.34%  '

$ws.Range("E11").Value = '  +6.23%  '

$ws.Range("D12").Value = '0.0933'
$ws.Range("E12").Value = '  +3.18%  '

$ws.Range("D13").Value = '2.089.91'
$ws.Range("E13").Value = '  +9.06%  '

$ws.Range("D14").Value = '1.833.03'
$ws.Range("E14").Value = '  +9.49%  '

$ws.Range("E15").Value = '  +8.42%  '

$ws.Range("D16").Value = '34.390.82'
$ws.Range("E16").Value = '  +12.16%  '

$ws.Range("D17").Value = '10.33'
$ws.Range("E17").Value = '  +4.69%  '

$ws.Range("D18").Value = '4.32'
$ws.Range("E18").Value = '  +7.55%  '

$ws.Range("D19").Value = '70.55'
$ws.Range("E19").Value = '  +6.56%  '

$ws.Range("D20").Value = '258.54'
$ws.Range("E20").Value = '  +6.92%  '

$ws.Range("D21").Value = '0.0₃0760'
$ws.Range("E21").Value = '  +5.54%  '

$ws.Range("D22").Value = '''1.00'
$ws.Range("E22").Value = '  +0.06%  '

$ws.Range("D23").Value = '10.66'
$ws.Range("E23").Value = '  +7.02%  '

$ws.Range("D24").Value = '4.35'
$ws.Range("E24").Value = '  +2.66%  '

$ws.Range("E25").Value = '  +3.92%  '

$ws.Range("D26").Value = '159.61'
$ws.Range("E26").Value = '  +0.34%  '

$ws.Range("D27").Value = '16.83'
$ws.Range("E27").Value = '  +6.51%  '

$ws.Range("E28").Value = '  +5.21%  '

$ws.Range("E29").Value = '  +7.73%  '

$ws.Range("E30").Value = '  -0.28%  '

$ws.Range("D31").Value = '''3.90'
$ws.Range("E31").Value = '  +12.88%  '

$ws.Range("D32").Value = '0.0526'
$ws.Range("E32").Value = '  +6.77%  '

$ws.Range("D33").Value = '1.22'
$ws.Range("E33").Value = '  +6.18%  '

$ws.Range("E34").Value = '  +8.04%  '

$ws.Range("D35").Value = '1.545.38'
$ws.Range("E35").Value = '  +3.14%  '

$ws.Range("D36").Value = '1.81'
$ws.Range("E36").Value = '  +2.10%  '

$ws.Range("D37").Value = '1.08'
$ws.Range("E37").Value = '  +5.29%  '

$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.0192'
$ws.Range("E38").Value = '  +7.97%  '

$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").Value = '0.639'
$ws.Range("E39").Value = '  +7.22%  '

$ws.Range("D40").Value = '84.87'
$ws.Range("E40").Value = '  +1.70%  '

$ws.Range("E41").Value = '  +5.29%  '

$ws.Range("E42").Value = '  +2.67%  '

$ws.Range("D43").Value = '0.916'
$ws.Range("E43").Value = '  +9.46%  '

$ws.Range("E44").Value = '  +6.32%  '

$ws.Range("D45").Value = '0.0528'
$ws.Range("E45").Value = '  +5.68%  '

$ws.Range("E46").Value = '  +6.15%  '

$ws.Range("D47").Value = '1.976.48'
$ws.Range("E47").Value = '  +9.05%  '

$ws.Range("D48").Value = '5.83'

$ws.Range("D49").Value = '12.23'
$ws.Range("E49").Value = '  +18.98%  '

$ws.Range("E50").Value = '  -0.06%  '

$ws.Range("D51").Value = '51.95'
$ws.Range("E51").Value = '  +4.13%  '
